$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A27").Value = "Command line for cypress"
$ws.Range("A29").Value = "npx cypress run --spec cypress/integration/examples/test9.js --headed --record --key 06499df6-3ffe-4a1d-872f-f3aa072ec3b4"
$ws.Range("A28").Value = "Run for cypress cloud and a specific test case in the spec and headed mode"

$ws.Range("A28").Select()
$excel.ActiveWindow.ScrollRow = 9
